$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update first three rows with new values
$ws.Range("A1").Value = 101089
$ws.Range("B1").Value = "Guerra Fria_atualizado"

$ws.Range("A2").Value = 101090
$ws.Range("B2").Value = "Primeira Guerra Mundial_atualizado"

$ws.Range("A3").Value = 101091
$ws.Range("B3").Value = "Segunda Guerra Mundial_atualizado"

# Remove row 4 entirely (delete the whole row so data shifts up / dimension shrinks)
$ws.Range("A4:B4").Delete()

# Select B5 as the active cell, to match the diff's selection state
$ws.Range("B5").Select()
